# Golden-test refresh: the title placeholders (and a couple of picture
# captions) on these slides were originally split into one run per word
# (e.g. "Slide" / " " / "1" / " " / "(Content)"). Re-set each one's text
# so the run-per-word split collapses into a single run per paragraph,
# matching the updated golden fixture.

$p = $ppt.ActivePresentation

function Set-ShapeText {
    param($shape, $text)
    # Force a real rewrite (collapsing multi-run text into one run) even
    # when the final joined text equals the current joined text, by
    # routing through a distinct intermediate value first.
    $shape.TextFrame.TextRange.Text = "__tmp__"
    $shape.TextFrame.TextRange.Text = $text
}

Set-ShapeText $p.Slides.Item(1).Shapes.Item(1)  "Slide 1 (Content)"
Set-ShapeText $p.Slides.Item(2).Shapes.Item(1)  "Slide 2 (Content)"
Set-ShapeText $p.Slides.Item(3).Shapes.Item(1)  "Slide 3 (Content)"
Set-ShapeText $p.Slides.Item(4).Shapes.Item(1)  "Slide 4 (Content)"
Set-ShapeText $p.Slides.Item(5).Shapes.Item(1)  "Slide 5 (Two Content)"
Set-ShapeText $p.Slides.Item(6).Shapes.Item(1)  "Slide 6 (Two Content Right)"
Set-ShapeText $p.Slides.Item(7).Shapes.Item(1)  "Slide 7 (Content with Caption)"
Set-ShapeText $p.Slides.Item(8).Shapes.Item(1)  "Slide 8 (Comparison)"
Set-ShapeText $p.Slides.Item(9).Shapes.Item(1)  "Slide 10 (Content)"
Set-ShapeText $p.Slides.Item(10).Shapes.Item(1) "Slide 11 (Content)"
Set-ShapeText $p.Slides.Item(11).Shapes.Item(1) "Slide 12 (Content)"

# Caption textboxes under the pictures.
Set-ShapeText $p.Slides.Item(6).Shapes.Item(3) "an image"
Set-ShapeText $p.Slides.Item(7).Shapes.Item(4) "An image"
Set-ShapeText $p.Slides.Item(8).Shapes.Item(4) "An image"
